# "naprawa błędów i hostowanie na serwerze"
#
# Adds a new log entry for Marek Wójcik (the first person block, columns
# B/C/D = Data/Plik/Linie) on row 37: the date moves from 2025-05-13
# (45790) to 2025-05-21 (45798), a new file "host aplikacji na serwerze"
# is logged with 300 lines added. Row 38's now-duplicate date is cleared.
# The D3/D4/G4/J4/M4/P4 formulas recalc automatically from the new D37
# input.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Arkusz1")

# Row 37: date, new file entry, line count
$ws.Range("B37").Value = 45798
$ws.Range("C37").Value = "host aplikacji na serwerze"
$ws.Range("D37").Value = 300

# Row 38: date cleared (was a duplicate of row 37's old date)
$ws.Range("B38").ClearContents()

# View state: scroll down a bit and move the selection, as in the saved
# workbook view.
$win = $excel.ActiveWindow
$win.ScrollRow = 5
$win.ScrollColumn = 1
$ws.Range("F40").Select()
